# Update cryptocurrency price/volume data to reflect latest snapshot
# (values refreshed by the scheduled GitHub Actions job)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.918.79"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "1.757.07"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.57%  "
$ws.Range("D5").Value = "'336.47"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").Value = "'0.3847"
$ws.Range("E7").Value = "  -1.47%  "
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("D9").Value = "'44.76"
$ws.Range("E9").Value = "  -6.44%  "
$ws.Range("D10").Value = "'1.110"
$ws.Range("E10").Value = "  -4.54%  "
$ws.Range("D11").Value = "'0.07201"
$ws.Range("E11").Value = "  -4.66%  "
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "'22.23"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("D14").Value = "'6.135"
$ws.Range("E14").Value = "  -5.13%  "
$ws.Range("D15").Value = "'7.143"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").Value = "1.757.78"
$ws.Range("E16").Value = "  -1.43%  "
$ws.Range("D17").Value = "'0.00001055"
$ws.Range("E17").Value = "  -3.38%  "
$ws.Range("D18").Value = "'0.06609"
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("D19").Value = "'79.15"
$ws.Range("E19").Value = "  -5.76%  "
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").Value = "'6.211"
$ws.Range("E21").Value = "  -5.28%  "
$ws.Range("D22").Value = "'16.57"
$ws.Range("D23").Value = "27.938.14"
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("D24").Value = "'11.60"
$ws.Range("E24").Value = "  -6.39%  "
$ws.Range("D25").Value = "'2.385"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("D26").Value = "'152.18"
$ws.Range("E26").Value = "  -2.27%  "
$ws.Range("D27").Value = "'19.78"
$ws.Range("E27").Value = "  -6.70%  "
$ws.Range("D28").Value = "'2.295"
$ws.Range("E28").Value = "  -8.73%  "
$ws.Range("D29").Value = "1.959.63"
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("D30").Value = "'1.272"
$ws.Range("E30").Value = "  -15.63%  "
$ws.Range("D31").Value = "'131.81"
$ws.Range("E31").Value = "  -4.55%  "
$ws.Range("D32").Value = "'4.023"
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("D33").Value = "'5.810"
$ws.Range("E33").Value = "  -7.67%  "
$ws.Range("D34").Value = "'0.08819"
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("D35").Value = "'12.15"
$ws.Range("E35").Value = "  -7.10%  "
$ws.Range("D36").Value = "'0.6587"
$ws.Range("E36").Value = "  -4.98%  "
$ws.Range("D37").Value = "'0.06180"
$ws.Range("E37").Value = "  -4.77%  "
$ws.Range("D38").Value = "'0.02284"
$ws.Range("E38").Value = "  -7.79%  "
$ws.Range("B39").Value = "WEMIXTOKEN"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").Value = "'1.521"
$ws.Range("E39").Value = "  -4.09%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "'5.127"
$ws.Range("E40").Value = "  -7.20%  "
$ws.Range("D41").Value = "'0.2104"
$ws.Range("E41").Value = "  -6.05%  "
$ws.Range("D42").Value = "'1.206"
$ws.Range("E42").Value = "  -4.74%  "
$ws.Range("D43").Value = "'7.961"
$ws.Range("E43").Value = "  -6.90%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("D45").Value = "'13.80"
$ws.Range("E45").Value = "  -5.85%  "
$ws.Range("D46").Value = "'3.827"
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("D47").Value = "'0.6013"
$ws.Range("E47").Value = "  -5.94%  "
$ws.Range("D48").Value = "'126.05"
$ws.Range("E48").Value = "  -5.92%  "
$ws.Range("E49").Value = "  -7.19%  "
$ws.Range("E50").Value = "  +3.94%  "
$ws.Range("D51").Value = "'1.168"
$ws.Range("E51").Value = "  +0.00%  "
